# Generate Report for Handback
#
# The handback process has produced target/handback files that are in sync
# with en-US, so for each tracked document we now know the "Latest Target
# File" (F) and "Latest Handback File" (G). The Status column is updated to
# reflect the handback, and the Latest Handback DateTime is stamped once the
# handback has actually completed (true so far only for de-de).

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")
$wsOv = $wb.Worksheets.Item("Overview")

# ---------------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# Replace across every sheet that shows this status so the shared string is
# updated everywhere it is used (Overview roll-up included).
# ---------------------------------------------------------------------------
$wsOv.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
$wsZh.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
$wsDe.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")

# ---------------------------------------------------------------------------
# Latest Handback DateTime placeholder "0001-01-01 00:00:00" becomes a real
# timestamp. Update both sheets first (so the un-handled-back placeholder is
# replaced consistently), then stamp de-de with its own, later, completion
# time since that handback has now actually finished.
# ---------------------------------------------------------------------------
$wsZh.Cells.Replace("0001-01-01 00:00:00", "2016-03-20 08:49:38")
$wsDe.Cells.Replace("0001-01-01 00:00:00", "2016-03-20 08:49:38")

# ---------------------------------------------------------------------------
# zh-cn sheet: populate Latest Target File (F) / Latest Handback File (G)
# ---------------------------------------------------------------------------

# Row 2 : 109dcc50-6bd8-4e8b-b586-bcfafe5be738
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md",
    "",
    "",
    "109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a53dabaacb2427959445960a82da41bd3a07e78/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/109dcc50-6bd8-4e8b-b586-bcfafe5be738.15ebeab342c42b42b91b326c5aeb3f11342164bb.zh-cn.xlf",
    "",
    "",
    "109dcc50-6bd8-4e8b-b586-bcfafe5be738.15ebeab342c42b42b91b326c5aeb3f11342164bb.zh-cn.xlf"
)

# Row 3 : a3dae023-59cc-481c-b11e-951bd45f84ca
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md",
    "",
    "",
    "a3dae023-59cc-481c-b11e-951bd45f84ca.md"
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4a53dabaacb2427959445960a82da41bd3a07e78/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a3dae023-59cc-481c-b11e-951bd45f84ca.b46bf5b7f4eb4de4fdd58708686169a2bef7598a.zh-cn.xlf",
    "",
    "",
    "a3dae023-59cc-481c-b11e-951bd45f84ca.b46bf5b7f4eb4de4fdd58708686169a2bef7598a.zh-cn.xlf"
)

# ---------------------------------------------------------------------------
# de-de sheet: populate Latest Target File (F) / Latest Handback File (G)
# ---------------------------------------------------------------------------

# Row 2 : 109dcc50-6bd8-4e8b-b586-bcfafe5be738
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/109dcc50-6bd8-4e8b-b586-bcfafe5be738.md",
    "",
    "",
    "109dcc50-6bd8-4e8b-b586-bcfafe5be738.md"
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64650fa90b9746eaa4376300910808b45395f0de/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/109dcc50-6bd8-4e8b-b586-bcfafe5be738.15ebeab342c42b42b91b326c5aeb3f11342164bb.de-de.xlf",
    "",
    "",
    "109dcc50-6bd8-4e8b-b586-bcfafe5be738.15ebeab342c42b42b91b326c5aeb3f11342164bb.de-de.xlf"
)

# Latest Handback DateTime (H2) is now known for this file: the de-de
# handback actually completed, so give it its own, distinct timestamp.
$wsDe.Range("H2").Value = "2016-03-20 08:49:43"

# Row 3 : a3dae023-59cc-481c-b11e-951bd45f84ca
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/332d22984f586bb1cf73b6fb171475f982528f71/e2e/a3dae023-59cc-481c-b11e-951bd45f84ca.md",
    "",
    "",
    "a3dae023-59cc-481c-b11e-951bd45f84ca.md"
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64650fa90b9746eaa4376300910808b45395f0de/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a3dae023-59cc-481c-b11e-951bd45f84ca.b46bf5b7f4eb4de4fdd58708686169a2bef7598a.de-de.xlf",
    "",
    "",
    "a3dae023-59cc-481c-b11e-951bd45f84ca.b46bf5b7f4eb4de4fdd58708686169a2bef7598a.de-de.xlf"
)

# Latest Handback DateTime (H3) is now known for this file
$wsDe.Range("H3").Value = "2016-03-20 08:49:43"
